$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(22)
$tr = $shp.TextFrame.TextRange

# Paragraph 2: merge split runs ("We " + "trained..." -> one run;
# "11583 rows..." + "81.24%" + " " + "accuracy..." -> one run)
$para2 = $tr.Paragraphs(2)
$c1 = $para2.Characters(1, 121)
$c1.Text = 'We trained our models on the English dataset. By just keeping gender classes “male” and “female” we could get a balanced '
$c2 = $para2.Characters(135, 115)
$c2.Text = '11583 rows. This gave us accuracy of 91.80% using the Random Forest algorithm and 81.24% accuracy using linear SVM.'

# Paragraph 5: merge "Estonian...86.22" into one run, leaving the
# trailing "%" as its own separate run.
$para5 = $tr.Paragraphs(5)
$c5 = $para5.Characters(1, 60)
$c5.Text = 'Estonian – RF accuracy: 87.92% , SVM(Linear) accuracy: 86.22'

# Append the new paragraphs (plain text first; run-split below).
$tr.InsertAfter([char]13 + 'Swedish – RF accuracy: 89.72% , SVM(Linear) accuracy: 84.64%' + [char]13 + 'Russian – RF accuracy: 90.78% , SVM(Linear) accuracy: 87.24%' + [char]13 + 'Chinese – RF accuracy: 91.03% , SVM(Linear) accuracy: 86.98%' + [char]13 + 'Italian – RF accuracy: 80.51% , SVM(Linear) accuracy: 85.43%' + [char]13 + [char]13 + [char]13 + [char]13) | Out-Null

# Re-split each new paragraph into its target run boundaries by
# re-asserting each run's text over its own character span (this
# creates a run boundary without altering the visible text).
$para = $tr.Paragraphs(6)
$para.Characters(1, 8).Text = 'Swedish '
$para.Characters(9, 15).Text = '– RF accuracy: '
$para.Characters(24, 7).Text = '89.72% '
$para.Characters(31, 24).Text = ', SVM(Linear) accuracy: '
$para.Characters(55, 6).Text = '84.64%'

$para = $tr.Paragraphs(7)
$para.Characters(1, 8).Text = 'Russian '
$para.Characters(9, 15).Text = '– RF accuracy: '
$para.Characters(24, 7).Text = '90.78% '
$para.Characters(31, 24).Text = ', SVM(Linear) accuracy: '
$para.Characters(55, 6).Text = '87.24%'

$para = $tr.Paragraphs(8)
$para.Characters(1, 10).Text = 'Chinese – '
$para.Characters(11, 11).Text = 'RF accuracy'
$para.Characters(22, 9).Text = ': 91.03% '
$para.Characters(31, 24).Text = ', SVM(Linear) accuracy: '
$para.Characters(55, 6).Text = '86.98%'

$para = $tr.Paragraphs(9)
$para.Characters(1, 8).Text = 'Italian '
$para.Characters(9, 15).Text = '– RF accuracy: '
$para.Characters(24, 7).Text = '80.51% '
$para.Characters(31, 24).Text = ', SVM(Linear) accuracy: '
$para.Characters(55, 6).Text = '85.43%'

